$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8 (shifts rows 8..24 down to 9..25),
# copying formatting from the row above (row 7, "directors") so that
# the new row 8 looks the same as row 7 then gets its text updated to
# "additional_directors".
$ws.Rows.Item(7).Copy() | Out-Null
$ws.Rows.Item(8).Insert(-4161)   # xlShiftDown

# Update the two cells that differ from the copied "directors" row.
$ws.Cells.Item(8, 1).Value = "additional_directors"
$ws.Cells.Item(8, 7).Value = "additional_directors"

# Update selection to match the target workbook state.
$ws.Range("G9").Select() | Out-Null

$wb.Save()
